$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg) for rows 2-21.
# Values derived from the updated weekly dataset (rows reshuffled with new dates/prices/origins).
$data = @(
    @{Row=2; Year=2021; Month=2; Day=9; J=100; K=25000; L=26000; M=25500; O='Región del Maule'; P=1020}
    @{Row=3; Year=2021; Month=6; Day=2; J=60; K=30000; L=32000; M=31000; O='Región Metropolitana'; P=1240}
    @{Row=4; Year=2021; Month=2; Day=3; J=100; K=35000; L=36000; M=35500; O='Región del Maule'; P=1420}
    @{Row=5; Year=2021; Month=3; Day=10; J=100; K=22000; L=24000; M=23000; O='Región del Maule'; P=920}
    @{Row=6; Year=2021; Month=1; Day=7; J=100; K=25000; L=26000; M=25500; O='Región de O''Higgins'; P=1020}
    @{Row=7; Year=2021; Month=5; Day=4; J=100; K=26000; L=28000; M=27000; O='Región del Maule'; P=1080}
    @{Row=8; Year=2021; Month=4; Day=27; J=100; K=30000; L=32000; M=31000; O='Región Metropolitana'; P=1240}
    @{Row=9; Year=2021; Month=2; Day=24; J=100; K=27000; L=28000; M=27500; O='Región Metropolitana'; P=1100}
    @{Row=10; Year=2021; Month=1; Day=14; J=100; K=32000; L=34000; M=33000; O='Región del Maule'; P=1320}
    @{Row=11; Year=2020; Month=12; Day=29; J=100; K=30000; L=32000; M=31000; O='Región Metropolitana'; P=1240}
    @{Row=12; Year=2020; Month=12; Day=23; J=100; K=42000; L=44000; M=43000; O='Región de O''Higgins'; P=1720}
    @{Row=13; Year=2021; Month=1; Day=20; J=100; K=26000; L=28000; M=27000; O='Región del Maule'; P=1080}
    @{Row=14; Year=2021; Month=3; Day=17; J=100; K=22000; L=24000; M=23000; O='Región del Maule'; P=920}
    @{Row=15; Year=2021; Month=4; Day=8; J=100; K=20000; L=22000; M=21000; O='Región del Maule'; P=840}
    @{Row=16; Year=2021; Month=4; Day=22; J=100; K=28000; L=30000; M=29000; O='Región del Maule'; P=1160}
    @{Row=17; Year=2021; Month=4; Day=30; J=100; K=26000; L=27000; M=26500; O='Región Metropolitana'; P=1060}
    @{Row=18; Year=2021; Month=3; Day=24; J=100; K=28000; L=30000; M=29000; O='Región del Maule'; P=1160}
    @{Row=19; Year=2021; Month=5; Day=26; J=100; K=28000; L=30000; M=29000; O='Región Metropolitana'; P=1160}
    @{Row=20; Year=2021; Month=5; Day=12; J=100; K=32000; L=34000; M=33000; O='Región Metropolitana'; P=1320}
    @{Row=21; Year=2021; Month=2; Day=17; J=100; K=25000; L=26000; M=25500; O='Región del Maule'; P=1020}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = (Get-Date -Year $item.Year -Month $item.Month -Day $item.Day -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
}
